# Update "Pais" sheet with refreshed COVID-19 country stats and provincias
# de España data (per commit "Update countries & provincias Spain").
#
# The source feed re-ranked a handful of countries (ties / near-ties in
# "Casos totales") causing their rows to swap, and refreshed the numbers
# for several countries. The "last updated" timestamp is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados ..." banner -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 22:55"

# --- Helper: write a full data row (Pais, Casos totales, Nuevos casos, --
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) -------
function Set-CountryRow($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Range("A$row").Value = $pais
    $ws.Range("B$row").Value = $casosTotales
    $ws.Range("C$row").Value = $nuevosCasos
    $ws.Range("D$row").Value = $casosActivos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $casosCriticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Rows whose underlying counts were refreshed (country stays in place)
Set-CountryRow 4   "Estados Unidos" 2350366 19788 977542 1250609 0 235 122215
Set-CountryRow 10  "Peru"           254936  3598  141967 104924  0 184 8045
Set-CountryRow 14  "Alemania"       191575  359   174900 7713    0 1   8962
Set-CountryRow 18  "Francia"        160377  284   74372  56365   0 7   29640
Set-CountryRow 50  "Barein"         21331   0     15790  5478    0 3   63
Set-CountryRow 51  "Israel"         20778   145   15694  4778    0 1   306
Set-CountryRow 82  "Guinea"         4988    28    3669   1292    0 0   27
Set-CountryRow 126 "Niger"          1036    1     911    58      0 0   67

# Rows that swapped rank with their neighbour (updated numbers + new order)
Set-CountryRow 73  "Costa de Marfil" 7492 216 3068 4370 0 2 54
Set-CountryRow 74  "Australia"       7461 25  6896 463  0 0 102

# Tied rows that swapped order (values identical, only the country changes)
Set-CountryRow 202 "Fiyi"             18 0 18 0 0 0 0
Set-CountryRow 203 "Dominica"         18 0 18 0 0 0 0
Set-CountryRow 207 "Groenlandia"      13 0 13 0 0 0 0
Set-CountryRow 208 "Islas Malvinas"   13 0 13 0 0 0 0
Set-CountryRow 213 "Papua Nueva Guinea"      8 0 8 0 0 0 0
Set-CountryRow 214 "Islas Virgenes Britanicas" 8 0 7 0 0 0 1
